$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3735.6924
$ws.Range("I64").Value = 3418.6667
$ws.Range("J64").Value = 4105.5557
$ws.Range("K64").Value = 3418.6667
$ws.Range("L64").Value = 4105.5557
$ws.Range("M64").Value = -3170.6667
$ws.Range("N64").Value = -4601.5557

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3735.6924
$ws.Range("I67").Value = 3418.6667
$ws.Range("J67").Value = 4105.5557
$ws.Range("K67").Value = 3418.6667
$ws.Range("L67").Value = 4105.5557
$ws.Range("M67").Value = -2560.6667
$ws.Range("N67").Value = -5821.5557

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 527.5714
$ws.Range("I80").Value = 580.6
$ws.Range("J80").Value = 498.1111
$ws.Range("K80").Value = 1741.8
$ws.Range("L80").Value = 1494.3333
$ws.Range("M80").Value = -743.8000000000002
$ws.Range("N80").Value = -3490.3333

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 527.5714
$ws.Range("I83").Value = 580.6
$ws.Range("J83").Value = 498.1111
$ws.Range("K83").Value = 5225.400000000001
$ws.Range("L83").Value = 4482.9999
$ws.Range("M83").Value = -233.4000000000005
$ws.Range("N83").Value = -14466.9999

# ALC row 105
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 25000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 25000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1206.3158
$ws.Range("I107").Value = 1163.3334
$ws.Range("J107").Value = 1980
$ws.Range("K107").Value = 1163.3334
$ws.Range("L107").Value = 1980
$ws.Range("M107").Value = 756.6666
$ws.Range("N107").Value = -5820

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1816.5312
$ws.Range("I113").Value = 1720.9
$ws.Range("J113").Value = 1860
$ws.Range("K113").Value = 1720.9
$ws.Range("L113").Value = 1860
$ws.Range("M113").Value = 1533.1
$ws.Range("N113").Value = -8368

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 48780
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 48780
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 48780
$ws.Range("N123").Value = -58580

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5706477.5
$ws.Range("I137").Value = 8622012
$ws.Range("J137").Value = 69777.60000000001
$ws.Range("K137").Value = 25866036
$ws.Range("L137").Value = 209332.8
$ws.Range("M137").Value = -25863486
$ws.Range("N137").Value = -214432.8

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18874790
$ws.Range("I32").Value = 22224508
$ws.Range("J32").Value = 32624.75
$ws.Range("K32").Value = 22224508
$ws.Range("L32").Value = 32624.75
$ws.Range("M32").Value = -22224221
$ws.Range("N32").Value = -33198.75

# BSM row 36
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 833.1667
$ws.Range("I36").Value = 833.1667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 833.1667
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -299.1667

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 274.7647
$ws.Range("I22").Value = 274.7647
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 274.7647
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 75.2353
$ws.Range("N22").ClearContents()

# CRP row 80
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 17705.334
$ws.Range("I80").Value = 17116
$ws.Range("J80").Value = 18000
$ws.Range("K80").Value = 17116
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = -15993
$ws.Range("N80").Value = -20246

# CRP row 83
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 17705.334
$ws.Range("I83").Value = 17116
$ws.Range("J83").Value = 18000
$ws.Range("K83").Value = 51348
$ws.Range("L83").Value = 54000
$ws.Range("M83").Value = -45732
$ws.Range("N83").Value = -65232

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 981.619
$ws.Range("I122").Value = 963.8946999999999
$ws.Range("J122").Value = 1150
$ws.Range("K122").Value = 2891.6841
$ws.Range("L122").Value = 3450
$ws.Range("M122").Value = -441.6840999999999
$ws.Range("N122").Value = -8350

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 40864.1
$ws.Range("I34").Value = 64
$ws.Range("J34").Value = 48064.117
$ws.Range("K34").Value = 192
$ws.Range("L34").Value = 144192.351
$ws.Range("M34").Value = -108
$ws.Range("N34").Value = -144360.351

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1742
$ws.Range("I39").Value = 355
$ws.Range("J39").Value = 2666.6667
$ws.Range("K39").Value = 1065
$ws.Range("L39").Value = 8000.000100000001
$ws.Range("M39").Value = -771
$ws.Range("N39").Value = -8588.000100000001

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2678.6875
$ws.Range("I55").Value = 567
$ws.Range("J55").Value = 3269.96
$ws.Range("K55").Value = 1701
$ws.Range("L55").Value = 9809.880000000001
$ws.Range("M55").Value = -1524
$ws.Range("N55").Value = -10163.88

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 826.7406999999999
$ws.Range("I68").Value = 676.5417
$ws.Range("J68").Value = 1045.2122
$ws.Range("K68").Value = 2029.6251
$ws.Range("L68").Value = 3135.6366
$ws.Range("M68").Value = -1218.6251
$ws.Range("N68").Value = -4757.6366

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 826.7406999999999
$ws.Range("I71").Value = 676.5417
$ws.Range("J71").Value = 1045.2122
$ws.Range("K71").Value = 6088.8753
$ws.Range("L71").Value = 9406.909799999999
$ws.Range("M71").Value = -2032.8753
$ws.Range("N71").Value = -17518.9098

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 956.63635
$ws.Range("I131").Value = 557.25
$ws.Range("J131").Value = 1011.7241
$ws.Range("K131").Value = 1671.75
$ws.Range("L131").Value = 3035.1723
$ws.Range("M131").Value = 3368.25
$ws.Range("N131").Value = -13115.1723

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7688.1055
$ws.Range("I113").Value = 1804.6923
$ws.Range("J113").Value = 20435.5
$ws.Range("K113").Value = 1804.6923
$ws.Range("L113").Value = 20435.5
$ws.Range("M113").Value = 365.3077000000001
$ws.Range("N113").Value = -24775.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 166668450
$ws.Range("I46").Value = 250000930
$ws.Range("J46").Value = 3490
$ws.Range("K46").Value = 250000930
$ws.Range("L46").Value = 3490
$ws.Range("M46").Value = -250000742
$ws.Range("N46").Value = -3866

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6595.722
$ws.Range("I93").Value = 8252.071
$ws.Range("J93").Value = 798.5
$ws.Range("K93").Value = 8252.071
$ws.Range("L93").Value = 798.5
$ws.Range("M93").Value = -7004.071
$ws.Range("N93").Value = -3294.5

# WVR row 24
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 4000
$ws.Range("I24").Value = 4000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -3770
$ws.Range("N24").ClearContents()

# WVR row 29
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1673000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1673000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1673000
$ws.Range("N29").Value = -1673580

# WVR row 75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 12096.5
$ws.Range("I75").Value = 10063
$ws.Range("J75").Value = 14130
$ws.Range("K75").Value = 10063
$ws.Range("L75").Value = 14130
$ws.Range("M75").Value = -9127
$ws.Range("N75").Value = -16002

# WVR row 78
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 12096.5
$ws.Range("I78").Value = 10063
$ws.Range("J78").Value = 14130
$ws.Range("K78").Value = 30189
$ws.Range("L78").Value = 42390
$ws.Range("M78").Value = -25509
$ws.Range("N78").Value = -51750

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1823.1052
$ws.Range("I132").Value = 1358.5
$ws.Range("J132").Value = 4301
$ws.Range("K132").Value = 4075.5
$ws.Range("L132").Value = 12903
$ws.Range("M132").Value = -1545.5
$ws.Range("N132").Value = -17963
